$d = $word.ActiveDocument
$wns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Paragraph right after "1. Introduzione" (currently an empty Titolo2
#    paragraph) becomes a normal paragraph containing the placeholder text
#    "/*INSERIRE INTRODUZIONE*/" (heading style removed).
# ---------------------------------------------------------------------------
$pIntro = $d.Paragraphs(117)
$pIntro.Range.InsertXML("<w:p$wns><w:r><w:t>/*INSERIRE INTRODUZIONE*/</w:t></w:r></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 2) Paragraph "2.1 Panoramica" (Titolo2 heading) right after
#    "2. Sistema proposto" is cleared out and turned into two blank
#    paragraphs (heading + text removed).
# ---------------------------------------------------------------------------
$pPanoramica = $d.Paragraphs(120)
$pPanoramica.Range.InsertXML("<w:p$wns></w:p><w:p$wns></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 3) Log-out use case ("RF2.2.7"): remove the stray <w:lastRenderedPageBreak/>
#    from the "Questa funzionalita` permette di" run.
#    (index +1 vs. the original document because step 2 above turned one
#    paragraph into two)
# ---------------------------------------------------------------------------
$pLogout = $d.Paragraphs(155)
$pLogout.Range.InsertXML("<w:p$wns w:rsidR=`"00ED507F`" w:rsidRPr=`"00677F22`" w:rsidRDefault=`"00ED507F`" w:rsidP=`"00ED507F`"><w:pPr><w:pStyle w:val=`"Paragrafoelenco`"/><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>Questa funzionalità permette di</w:t></w:r><w:r w:rsidRPr=`"00677F22`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> disconnettersi dal sito</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>.</w:t></w:r></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 4) "Effettuare un ordine" use case ("RF2.2.8"): add
#    <w:lastRenderedPageBreak/> to the "Priorita`" run (it moved here from
#    the Log-out paragraph above).
# ---------------------------------------------------------------------------
$pPriorita1 = $d.Paragraphs(159)
$pPriorita1.Range.InsertXML("<w:p$wns w:rsidR=`"00ED507F`" w:rsidRPr=`"00677F22`" w:rsidRDefault=`"00ED507F`" w:rsidP=`"00ED507F`"><w:pPr><w:pStyle w:val=`"Paragrafoelenco`"/><w:spacing w:line=`"276`" w:lineRule=`"auto`"/><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r w:rsidRPr=`"00677F22`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Priorità</w:t></w:r><w:r w:rsidRPr=`"00677F22`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>: Alta</w:t></w:r></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 5) "RNF2.3.1.3" requirement: remove the stray <w:lastRenderedPageBreak/>
#    from the "Il sistema deve, in caso di input errato..." run.
# ---------------------------------------------------------------------------
$pRnf313 = $d.Paragraphs(193)
$pRnf313.Range.InsertXML("<w:p$wns w:rsidR=`"00ED507F`" w:rsidRPr=`"00677F22`" w:rsidRDefault=`"00ED507F`" w:rsidP=`"00ED507F`"><w:pPr><w:pStyle w:val=`"Paragrafoelenco`"/><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Il sistema deve, in caso di input errato da parte dell’utente durante la compilazione di un </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>form</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>, evidenziare i campi scorretti e far visualizzare un messaggio testuale che indichi come riempire correttamente il campo.</w:t></w:r></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 6) "RNF2.3.1.4" heading: add <w:lastRenderedPageBreak/> to the "RNF" run
#    (it moved here from the RNF2.3.1.3 paragraph above).
# ---------------------------------------------------------------------------
$pRnf314 = $d.Paragraphs(195)
$pRnf314.Range.InsertXML("<w:p$wns w:rsidR=`"00ED507F`" w:rsidRPr=`"00677F22`" w:rsidRDefault=`"00ED507F`" w:rsidP=`"00ED507F`"><w:pPr><w:pStyle w:val=`"Paragrafoelenco`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:lastRenderedPageBreak/><w:t>RNF</w:t></w:r><w:r w:rsidR=`"000D3EE6`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>2.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>3.1.4</w:t></w:r></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 7) "2.5 System Model" heading: remove the stray
#    <w:lastRenderedPageBreak/> from the "2.5 " run.
# ---------------------------------------------------------------------------
$pSystemModel = $d.Paragraphs(222)
$pSystemModel.Range.InsertXML("<w:p$wns w:rsidR=`"000D3EE6`" w:rsidRDefault=`"000D3EE6`" w:rsidP=`"000D3EE6`"><w:pPr><w:pStyle w:val=`"Titolo1`"/></w:pPr><w:r><w:t xml:space=`"preserve`">2.5 </w:t></w:r><w:r w:rsidR=`"00D07E5B`"><w:t>System Model</w:t></w:r></w:p>") | Out-Null

# ---------------------------------------------------------------------------
# 8) Remove 4 of the 5 consecutive empty paragraphs right before the
#    "Un utente registrato..." paragraph (registered-user feature list),
#    keeping the paragraph that already carries the text.
#    Must run AFTER all the paragraph-index-preserving edits above, since
#    it changes the total paragraph count.
#    (index +1 vs. the original document because step 2 above turned one
#    paragraph into two)
# ---------------------------------------------------------------------------
$pBlank1 = $d.Paragraphs(141)
$pBlank4 = $d.Paragraphs(144)
$rBlanks = $d.Range($pBlank1.Range.Start, $pBlank4.Range.End)
$rBlanks.Delete() | Out-Null
